# Implemented playtest1 gameflow and balancing
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7 (Tanuki): EVADE 10 -> 0 ---
$ws.Range("K7").Value = 0

# --- Row 9 (Ijiraq2): max_hp 40 -> 20 ---
$ws.Range("D9").Value = 20

# --- Insert a new row at 11 for "Ijiraq3" (pushes old rows 11-15 down to 12-16) ---
$ws.Rows.Item(11).Insert()

$ws.Range("A11").Value = "Ijiraq3"
$ws.Range("B11").Value = "spr_bt_ijiraq_placeholder"
$ws.Range("C11").Value = "enemy_general_1"
$ws.Range("D11").Value = 45
$ws.Range("E11").Value = 0
$ws.Range("F11").Value = 2
$ws.Range("G11").Value = 0.75
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0.5
$ws.Range("J11").Value = 1
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = "Attacker1"
$ws.Range("M11").Value = "none"
$ws.Range("N11").Value = 0
$ws.Range("O11").Value = 1.1000000000000001
$ws.Range("P11").Value = -1
$ws.Range("Q11").Value = 1
$ws.Range("R11").Value = "GROUP/DEFAULT"
$ws.Range("S11").Value = "spear "
$ws.Range("T11").Value = "ice"
$ws.Range("U11").Value = "null"
$ws.Range("V11").Value = "sword"
$ws.Range("W11").Value = "ice"
$ws.Range("X11").Value = "null"
$ws.Range("Y11").Value = "END"

# --- Selection moves to D11 ---
$ws.Range("D11").Select()
